$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text (avoids Excel auto-converting
# numeric-looking / date-looking strings into numbers or dates), while
# leaving the cell with the default ("Normal") style, i.e. no explicit
# style index in the saved XML - matching the unstyled data rows of the
# original sheet.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-RowData($row, $name, $aid, $species, $location, $intake) {
    Set-TextCell $row 1 $name
    Set-TextCell $row 2 $aid
    Set-TextCell $row 3 $species
    Set-TextCell $row 4 $location
    Set-TextCell $row 5 $intake
}

# New data set (replaces the previous animal list entirely, and adds
# five additional rows, rows 2-12).
Set-RowData 2  "PALM"     "58262248" "Cat"    "Feature Room 2"           "4/8/2025"
Set-RowData 3  "COLLIE"   "57884999" "Cat"    "If The Fur Fits"          "2/20/2025"
Set-RowData 4  "SPECK"    "52249653" "Dog"    "Dog Adoptions D"          "4/16/2025"
Set-RowData 5  "LAUREL"   "58289985" "Cat"    "Cat Adoption Condo Rooms" "4/11/2025"
Set-RowData 6  "Gyarados" "58096306" "Cat"    "Cat Treatment"            "3/28/2025"
Set-RowData 7  "Mochi"    "58353916" "Cat"    "Offsite Adoptions"        "4/22/2025"
Set-RowData 8  "SUDS"     "58598619" "Dog"    "Dog Adoptions A"          "5/29/2025"
Set-RowData 9  "Beau"     "58677023" "Dog"    "Dog Adoptions C"          "6/10/2025"
Set-RowData 10 "RAYNE"    "57710656" "Dog"    "If The Fur Fits"          "6/13/2025"
Set-RowData 11 "HEATH"    "58654173" "Dog"    "If The Fur Fits"          "6/6/2025"
Set-RowData 12 "MARINA"   "58706705" "Rabbit" "Adoptions Lobby"          "6/13/2025"

Write-Host "Updated sheet dimension:" $ws.UsedRange.Address()
